# Update "想去人数" (F column) counts on the 展览, 本地生活 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 69
$ws1.Range("F3").Value = 158
$ws1.Range("F6").Value = 2717
$ws1.Range("F9").Value = 7357
$ws1.Range("F11").Value = 7521
$ws1.Range("F14").Value = 5967
$ws1.Range("F15").Value = 3214
$ws1.Range("F16").Value = 3580
$ws1.Range("F20").Value = 430
$ws1.Range("F22").Value = 266
$ws1.Range("F23").Value = 2063
$ws1.Range("F28").Value = 942
$ws1.Range("F29").Value = 56
$ws1.Range("F30").Value = 2567
$ws1.Range("F31").Value = 1387
$ws1.Range("F32").Value = 3126
$ws1.Range("F33").Value = 131
$ws1.Range("F35").Value = 226
$ws1.Range("F37").Value = 459
$ws1.Range("F38").Value = 1205
$ws1.Range("F39").Value = 222

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 114

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 69
$ws4.Range("F5").Value = 158
$ws4.Range("F9").Value = 114
$ws4.Range("F10").Value = 2717
$ws4.Range("F14").Value = 7357
$ws4.Range("F16").Value = 7521
$ws4.Range("F18").Value = 5967
$ws4.Range("F19").Value = 3214
$ws4.Range("F20").Value = 3580
$ws4.Range("F24").Value = 430
$ws4.Range("F29").Value = 266
$ws4.Range("F30").Value = 2063
$ws4.Range("F36").Value = 942
$ws4.Range("F37").Value = 56
$ws4.Range("F38").Value = 2567
$ws4.Range("F39").Value = 1387
$ws4.Range("F41").Value = 3126
$ws4.Range("F42").Value = 131
$ws4.Range("F45").Value = 459
$ws4.Range("F46").Value = 1205
$ws4.Range("F47").Value = 222
